$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new "category" column before the existing "date" column (I),
# shifting date/legislator_name/legislator_id right by one (I->J, J->K, K->L).
$ws.Columns.Item(9).Insert()

# Header row values
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the header formatting (bold + border + centered/top alignment) used
# by the other header cells.
$ws.Range("H1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows values
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"
$ws.Range("I4").Value = "normal"

$ws.Range("M2").Value = "tmpabd41"
$ws.Range("M3").Value = "tmpabd41"
$ws.Range("M4").Value = "tmpabd41"

$ws.Range("N2").Value = 90
$ws.Range("N3").Value = 91
$ws.Range("N4").Value = 92

# Match the data-row formatting for the new trailing columns.
$ws.Range("H2").Copy()
$ws.Range("M2").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
